# "updating my reads and adding pkane goal data"
# Fill in the daily pages-read (F) and current-book (G) columns for the
# rows that were missing them (rows 77-95 of the reading-data sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("reading-data")

# row => (pages, book)
$dailyData = [ordered]@{
    "77" = @(6,3)
    "78" = @(0,3)
    "79" = @(0,3)
    "80" = @(2,3)
    "81" = @(6,3)
    "82" = @(0,3)
    "83" = @(0,3)
    "84" = @(0,3)
    "85" = @(0,3)
    "86" = @(5,3)
    "87" = @(0,3)
    "88" = @(0,3)
    "89" = @(5,3)
    "90" = @(0,3)
    "91" = @(8,3)
    "92" = @(0,3)
    "93" = @(0,3)
    "94" = @(0,3)
    "95" = @(5,3)
}

foreach ($key in $dailyData.Keys) {
    $row = [int]$key
    $pair = $dailyData[$key]
    $ws.Cells.Item($row, 6).Value = $pair[0]
    $ws.Cells.Item($row, 7).Value = $pair[1]
}

# Match the author's last on-screen selection/scroll position.
$ws.Activate() | Out-Null
$ws.Range("G76:G95").Select() | Out-Null
